# The workbook tracks transition-matrix data between reading-behaviour
# states. The state label "ScreenRecStarted" is being renamed to
# "0_unstated" everywhere it appears (the column/row header in G1, and
# the composite row labels in A27:A30 that are built from it).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "0_unstated"
$ws.Range("A27").Value = "0_unstated1_Scanning"
$ws.Range("A28").Value = "0_unstated3_Reading"
$ws.Range("A29").Value = "0_unstated5_Unknown "
$ws.Range("A30").Value = "0_unstated0_unstated"

# Move the active selection to G14, matching the sheet's last-saved
# selection state.
$ws.Range("G14").Select()
